# Update dates for the "Number" and "Testinomial" sections (and their
# "Responsive" sub-rows) on Sheet1.
#
# Rows affected:
#   29 (Number)       -> Finish Build (C29) was blank, set to 10/12/2021
#   32 (Responsive)    -> Start/Finish Build (B32/C32) were blank, set to 10/12/2021
#   39 (Testinomial)   -> Finish Build (C39) was blank, set to 10/12/2021
#   42 (Responsive)    -> Start/Finish Build (B42/C42) were blank, set to 10/12/2021

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = 44481   # 10/12/2021 as an Excel serial date

# --- Row 29 ("Number") ---
# C29 should pick up the same look as B29 (bold date style), then get its value.
$ws.Range("B29").Copy()
$ws.Range("C29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C29").Value = $newDate

# --- Row 32 ("Responsive" under Number) ---
# B32/C32 should pick up the same look as the analogous B30/C30 cells
# (non-bold date style used throughout the "Responsive" rows).
$ws.Range("B30").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("B32").Value = $newDate

$ws.Range("C30").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("C32").Value = $newDate

# --- Row 39 ("Testinomial") ---
$ws.Range("B39").Copy()
$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("C39").Value = $newDate

# --- Row 42 ("Responsive" under Testinomial) ---
$ws.Range("B40").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("B42").Value = $newDate

$ws.Range("C40").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("C42").Value = $newDate

$excel.CutCopyMode = $false

# Move the view so the edited area is visible/selected, matching where the
# author was working when they made this change.
$excel.Goto($ws.Range("A34"), $true)
$ws.Range("C39").Select()
